$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("observed_stns")

# Row 12 corresponds to station_id 36225 / station_code OSLEALN.
# The nve_id (C12) and dis_station_id (D12) values are no longer valid
# (the station switched to a modelled series), so clear those cells
# entirely (contents + formatting) and record an explanatory comment in
# column E instead.
$ws.Range("C12").Clear()
$ws.Range("D12").Clear()
$ws.Range("E12").Value = "Used to be 6.78.0 (dis_stn 626), but switched to modelled series due to difficulties obtaining data"

# Reflect the last worked-on cell, matching where the author left the
# selection after making this change.
$ws.Range("D12").Select()
